# Fruta / hortaliza, semanal
# Insert 4 new weekly price rows for Nectarín (Ruby Diamond x2, Venus x2)
# at the top of the "August Red / June Pearl" 2021-03-11 block, pushing the
# existing rows 190-203 down to 194-207.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before row 190 (formatting/style is inherited from
# the row above, same as a native Excel "Insert Copied Cells" / row insert).
$ws.Range("A190:A193").EntireRow.Insert()

# Common columns shared by every row in this block.
$mercadoId   = 8
$mercado     = "Terminal La Palmera de La Serena"
$region      = "Coquimbo"
$codreg      = 4
$tipo        = "Fruta"
$productoId  = 100103
$producto    = "Frutos de hueso (carozo)"
$categoriaId = 100103006
$categoria   = "Nectarín"
$origen      = "Región de O'Higgins"
$unidadBins  = "$/bins (420 kilos)"
$fecha       = 44578

# Row 190: Ruby Diamond, Especial
$r = 190
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Ruby Diamond"
$ws.Cells.Item($r, 12).Value = "Especial"
$ws.Cells.Item($r, 13).Value = 20
$ws.Cells.Item($r, 14).Value = 385000
$ws.Cells.Item($r, 15).Value = 390000
$ws.Cells.Item($r, 16).Value = 387500
$ws.Cells.Item($r, 17).Value = $unidadBins
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 923
$ws.Cells.Item($r, 20).Value = 420

# Row 191: Ruby Diamond, Primera
$r = 191
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Ruby Diamond"
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 20
$ws.Cells.Item($r, 14).Value = 325000
$ws.Cells.Item($r, 15).Value = 330000
$ws.Cells.Item($r, 16).Value = 327500
$ws.Cells.Item($r, 17).Value = $unidadBins
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 780
$ws.Cells.Item($r, 20).Value = 420

# Row 192: Venus, Especial
$r = 192
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Venus"
$ws.Cells.Item($r, 12).Value = "Especial"
$ws.Cells.Item($r, 13).Value = 20
$ws.Cells.Item($r, 14).Value = 455000
$ws.Cells.Item($r, 15).Value = 460000
$ws.Cells.Item($r, 16).Value = 457500
$ws.Cells.Item($r, 17).Value = $unidadBins
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 1089
$ws.Cells.Item($r, 20).Value = 420

# Row 193: Venus, Primera
$r = 193
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Venus"
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 20
$ws.Cells.Item($r, 14).Value = 415000
$ws.Cells.Item($r, 15).Value = 420000
$ws.Cells.Item($r, 16).Value = 417500
$ws.Cells.Item($r, 17).Value = $unidadBins
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 994
$ws.Cells.Item($r, 20).Value = 420
